$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 35-42 keep the same species but their weight (G) is zeroed out
for ($r = 35; $r -le 42; $r++) {
    $ws.Cells.Item($r, 7).Value = 0
}

# Rows 43-58 take on the identity (E, F, G, H) of the OLD rows 44-61 (shift up by one,
# dropping the old row 43 "Eggs of Murex"), with G always becoming 0.
$shifted = @(
    @{E='Eggs of Raja sp';        F='EGGSRAJ'; H=36},
    @{E='Eggs of Tonna';          F='EGGSTON'; H=1},
    @{E='Galeodea echinophora';   F='GALEECH'; H=1},
    @{E='Goneplax rhomboides';    F='GONERHO'; H=12},
    @{E='Medorippe lanata';       F='MEDOLAN'; H=22},
    @{E='Modiolus barbatus';      F='MODIBAR'; H=1},
    @{E='Mytilus galloprovincialis'; F='MYTGALL'; H=1},
    @{E='Nucula sulcata';         F='NUCUSUL'; H=1},
    @{E='Ostrea edulis';          F='OSTREDU'; H=1},
    @{E='Pagurus anachoretus';    F='PAGUANA'; H=1},
    @{E='Parthenope angulifrons'; F='PARTANG'; H=1},
    @{E='Polycarpa sp';           F='POLYSP';  H=1},
    @{E='Solenocera membranacea'; F='SOLOMEM'; H=1},
    @{E='Trachythyone elongata';  F='TRACELO'; H=2},
    @{E='Trachythyone tergestina';F='TRACTER'; H=65},
    @{E='Turritella communis';    F='TURRCOM'; H=46}
)

$r = 43
foreach ($row in $shifted) {
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = 0
    $ws.Cells.Item($r, 8).Value = $row.H
    $r++
}

# Old rows 59-62 ("Trachythyone elongata"..."Wood NA") no longer exist; delete them entirely
$deleteRange = $ws.Range("A59:K62")
$deleteRange.Delete()
